$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -22.0825
$ws.Range("A14").Value = -21.82489999999999
$ws.Range("D15").Value = -8.551300000000003
$ws.Range("A16").Value = -22.01610000000001
$ws.Range("A21").Value = -19.55079999999998
$ws.Range("D21").Value = -7.9933
$ws.Range("D22").Value = -8.140200000000004
$ws.Range("A23").Value = -20.25319999999997
$ws.Range("D24").Value = -7.639499999999999
$ws.Range("A25").Value = -21.81869999999998
$ws.Range("A26").Value = -21.04649999999997
$ws.Range("D27").Value = -7.674100000000002
$ws.Range("D28").Value = -7.887900000000001
$ws.Range("A29").Value = -21.14179999999997
$ws.Range("D36").Value = -6.937000000000003
$ws.Range("D39").Value = -7.943599999999999
$ws.Range("A40").Value = -19.75879999999999
$ws.Range("D45").Value = -7.036100000000005
$ws.Range("D48").Value = -7.211299999999993
$ws.Range("D49").Value = -7.888100000000001
$ws.Range("D52").Value = -7.845500000000005
$ws.Range("A53").Value = -22.0105
$ws.Range("D53").Value = -7.763099999999997
$ws.Range("D54").Value = -7.992000000000006
$ws.Range("A57").Value = -22.77530000000002
$ws.Range("D57").Value = -8.320799999999993
$ws.Range("A59").Value = -22.43180000000001
$ws.Range("A65").Value = -21.87889999999998
$ws.Range("A69").Value = -21.57889999999999
$ws.Range("D70").Value = -6.889600000000001
$ws.Range("D71").Value = -6.860799999999998
$ws.Range("A79").Value = -20.33300000000001
$ws.Range("A83").Value = -21.7389
$ws.Range("D86").Value = -8.87990000000001
$ws.Range("D87").Value = -8.180899999999996
$ws.Range("D89").Value = -8.654400000000003
$ws.Range("A91").Value = -20.75419999999997
$ws.Range("A93").Value = -21.43770000000001
$ws.Range("A100").Value = -22.0022
$ws.Range("D101").Value = -8.1067
$ws.Range("A103").Value = -21.67299999999999
